$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook-level structure protection record
$wb.Protect()

# The sheet's default ("Normal") font changes from Arial to Calibri
$wb.Styles.Item("Normal").Font.Name = "Calibri"

# Two new "Checkers" rows of data
$ws.Range("A3").Value = "dfbn"
$ws.Range("B3").Value = "sfgn"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = "age 6-8"
$ws.Range("F3").Value = "Checkers"

$ws.Range("A4").Value = "srfgn"
$ws.Range("B4").Value = "wrh"
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = "age 6-8"
$ws.Range("F4").Value = "Checkers"

$ws.Range("A1").Select()
